$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 887
$ws.Range("I28").Value = 882.25
$ws.Range("K28").Value = 882.25
$ws.Range("M28").Value = -397.25

$ws.Range("H42").Value = 1470.4667
$ws.Range("I42").Value = 973.4
$ws.Range("J42").Value = 2464.6
$ws.Range("K42").Value = 2920.2
$ws.Range("L42").Value = 7393.799999999999
$ws.Range("M42").Value = -2690.2
$ws.Range("N42").Value = -7853.799999999999

$ws.Range("H100").Value = 3221.8518
$ws.Range("I100").Value = 4165.75
$ws.Range("J100").Value = 1848.909
$ws.Range("K100").Value = 4165.75
$ws.Range("L100").Value = 1848.909
$ws.Range("M100").Value = -3624.75
$ws.Range("N100").Value = -2930.909

$ws.Range("H111").Value = 1359.2174
$ws.Range("I111").Value = 1411.2307
$ws.Range("J111").Value = 1291.6
$ws.Range("K111").Value = 4233.6921
$ws.Range("L111").Value = 3874.8
$ws.Range("M111").Value = -1166.6921
$ws.Range("N111").Value = -10008.8

$ws.Range("H137").Value = 1500
$ws.Range("I137").Value = 1500
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 4500
$ws.Range("L137").ClearContents()
$ws.Range("M137").Value = -1950
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1102.0416
$ws.Range("I32").Value = 1102.0416
$ws.Range("K32").Value = 1102.0416
$ws.Range("M32").Value = -815.0416

$ws.Range("H61").Value = 4863.143
$ws.Range("I61").Value = 4838.4287
$ws.Range("K61").Value = 4838.4287
$ws.Range("M61").Value = -4626.4287

$ws.Range("H74").Value = 1598334.2
$ws.Range("I74").Value = 687849
$ws.Range("J74").Value = 13889886
$ws.Range("K74").Value = 687849
$ws.Range("L74").Value = 13889886
$ws.Range("M74").Value = -686975
$ws.Range("N74").Value = -13891634

$ws.Range("H77").Value = 1598334.2
$ws.Range("I77").Value = 687849
$ws.Range("J77").Value = 13889886
$ws.Range("K77").Value = 3439245
$ws.Range("L77").Value = 69449430
$ws.Range("M77").Value = -3434877
$ws.Range("N77").Value = -69458166

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H132").Value = 500000000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H136").Value = 4863.143
$ws.Range("I136").Value = 4838.4287
$ws.Range("K136").Value = 14515.2861
$ws.Range("M136").Value = -11965.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2200.5
$ws.Range("I94").Value = 2200.5
$ws.Range("K94").Value = 2200.5
$ws.Range("M94").Value = -1749.5

$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1100
$ws.Range("J99").Value = 600
$ws.Range("K99").Value = 1100
$ws.Range("L99").Value = 600
$ws.Range("M99").Value = 398
$ws.Range("N99").Value = -3596

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2893.0527
$ws.Range("J31").Value = 3123
$ws.Range("L31").Value = 3123
$ws.Range("N31").Value = -3713

$ws.Range("H34").Value = 2893.0527
$ws.Range("J34").Value = 3123
$ws.Range("L34").Value = 3123
$ws.Range("N34").Value = -3527

$ws.Range("H134").Value = 5885179.5
$ws.Range("I134").Value = 2753.7856
$ws.Range("J134").Value = 33336500
$ws.Range("K134").Value = 8261.356800000001
$ws.Range("L134").Value = 100009500
$ws.Range("M134").Value = -5726.356800000001
$ws.Range("N134").Value = -100014570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 357.4
$ws.Range("I10").Value = 396.25
$ws.Range("J10").Value = 202
$ws.Range("K10").Value = 1188.75
$ws.Range("L10").Value = 606
$ws.Range("M10").Value = -1049.75
$ws.Range("N10").Value = -884

$ws.Range("H131").Value = 739029.06
$ws.Range("J131").Value = 1071000.9
$ws.Range("L131").Value = 3213002.7
$ws.Range("N131").Value = -3223082.7

$ws.Range("H134").Value = 5132.857
$ws.Range("I134").Value = 5132.857
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 15398.571
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -10328.571
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9104.571
$ws.Range("I70").Value = 8399.348
$ws.Range("J70").Value = 10456.25
$ws.Range("K70").Value = 8399.348
$ws.Range("L70").Value = 10456.25
$ws.Range("M70").Value = -8129.348
$ws.Range("N70").Value = -10996.25

$ws.Range("H73").Value = 9104.571
$ws.Range("I73").Value = 8399.348
$ws.Range("J73").Value = 10456.25
$ws.Range("K73").Value = 8399.348
$ws.Range("L73").Value = 10456.25
$ws.Range("M73").Value = -7463.348
$ws.Range("N73").Value = -12328.25

$ws.Range("H113").Value = 1275.8334
$ws.Range("I113").Value = 945.1429
$ws.Range("J113").Value = 1738.8
$ws.Range("K113").Value = 945.1429
$ws.Range("L113").Value = 1738.8
$ws.Range("M113").Value = 1224.8571
$ws.Range("N113").Value = -6078.8

$ws.Range("H126").Value = 20006.5
$ws.Range("I126").Value = 29999
$ws.Range("J126").Value = 10014
$ws.Range("K126").Value = 89997
$ws.Range("L126").Value = 30042
$ws.Range("M126").Value = -87527
$ws.Range("N126").Value = -34982

$ws.Range("H128").Value = 26000
$ws.Range("J128").Value = 26000
$ws.Range("L128").Value = 26000
$ws.Range("N128").Value = -35960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 20198.092
$ws.Range("I93").Value = 1812.9166
$ws.Range("K93").Value = 1812.9166
$ws.Range("M93").Value = -564.9166

$ws.Range("H136").Value = 22730326
$ws.Range("I136").Value = 2715.4595
$ws.Range("K136").Value = 8146.3785
$ws.Range("M136").Value = -5596.3785

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1473.1
$ws.Range("I81").Value = 1610.2222
$ws.Range("J81").Value = 239
$ws.Range("K81").Value = 3220.4444
$ws.Range("L81").Value = 478
$ws.Range("M81").Value = -2159.4444
$ws.Range("N81").Value = -2600

$ws.Range("H84").Value = 1473.1
$ws.Range("I84").Value = 1610.2222
$ws.Range("J84").Value = 239
$ws.Range("K84").Value = 16102.222
$ws.Range("L84").Value = 2390
$ws.Range("M84").Value = -10798.222
$ws.Range("N84").Value = -12998

$ws.Range("H107").Value = 1014.53845
$ws.Range("I107").Value = 1022.125
$ws.Range("K107").Value = 3066.375
$ws.Range("M107").Value = -1146.375

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H136").Value = 1036.5333
$ws.Range("I136").Value = 896.2857
$ws.Range("K136").Value = 2688.8571
$ws.Range("M136").Value = -138.8571000000002
